$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058005996960946
$ws.Range("D2").Value = 1.056605989075327
$ws.Range("E2").Value = 1.063696170159247
$ws.Range("F2").Value = 1.073978272471214
$ws.Range("I2").Value = 1.045527381813204
$ws.Range("J2").Value = 1.062999464318567
$ws.Range("K2").Value = 1.059343239533129
$ws.Range("L2").Value = 1.066414111728983
$ws.Range("M2").Value = 1.076668702223863
$ws.Range("N2").Value = 1.064509045485642

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059740620401336
$ws.Range("D3").Value = 1.057949509371131
$ws.Range("E3").Value = 1.065267267848745
$ws.Range("F3").Value = 1.075653167078425
$ws.Range("I3").Value = 1.046034609677996
$ws.Range("J3").Value = 1.064383011354173
$ws.Range("K3").Value = 1.060499084107653
$ws.Range("L3").Value = 1.067798375033153
$ws.Range("M3").Value = 1.078158527983149
$ws.Range("N3").Value = 1.065894557316734

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060859979086272
$ws.Range("D4").Value = 1.058816079259723
$ws.Range("E4").Value = 1.066280829821064
$ws.Range("F4").Value = 1.076734191208752
$ws.Range("I4").Value = 1.046359929583471
$ws.Range("J4").Value = 1.06527489073295
$ws.Range("K4").Value = 1.061243681514274
$ws.Range("L4").Value = 1.068690554513163
$ws.Range("M4").Value = 1.079119327759909
$ws.Range("N4").Value = 1.066787703266529

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.061329840635215
$ws.Range("D5").Value = 1.059179731201898
$ws.Range("E5").Value = 1.066706216399681
$ws.Range("F5").Value = 1.077188010774104
$ws.Range("I5").Value = 1.04649600729005
$ws.Range("J5").Value = 1.065649042941242
$ws.Range("K5").Value = 1.061555927809105
$ws.Range("L5").Value = 1.069064794106413
$ws.Range("M5").Value = 1.079522490924167
$ws.Range("N5").Value = 1.067162386813872

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.061408690865975
$ws.Range("D6").Value = 1.059240751925261
$ws.Range("E6").Value = 1.066777599153093
$ws.Range("F6").Value = 1.077264171778844
$ws.Range("I6").Value = 1.04651881521247
$ws.Range("J6").Value = 1.065711818548989
$ws.Range("K6").Value = 1.06160830974929
$ws.Range("L6").Value = 1.069127582122278
$ws.Range("M6").Value = 1.079590139778003
$ws.Range("N6").Value = 1.067225251570183

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060866260194912
$ws.Range("D7").Value = 1.058820940949694
$ws.Range("E7").Value = 1.066286516647558
$ws.Range("F7").Value = 1.076740257675144
$ws.Range("I7").Value = 1.046361750551872
$ws.Range("J7").Value = 1.065279893275441
$ws.Range("K7").Value = 1.061247856824944
$ws.Range("L7").Value = 1.068695558375216
$ws.Range("M7").Value = 1.0791247178026
$ws.Range("N7").Value = 1.066792712913203

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058592863707716
$ws.Range("D8").Value = 1.057060619769805
$ws.Range("E8").Value = 1.064227768390716
$ws.Range("F8").Value = 1.074544887773909
$ws.Range("I8").Value = 1.045699403473031
$ws.Range("J8").Value = 1.063467745419
$ws.Range("K8").Value = 1.059734554468208
$ws.Range("L8").Value = 1.066882668923757
$ws.Range("M8").Value = 1.077172869870779
$ws.Range("N8").Value = 1.06497799159895

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054562728175009
$ws.Range("D9").Value = 1.053936943153364
$ws.Range("E9").Value = 1.060576065608785
$ws.Range("F9").Value = 1.070654708808215
$ws.Range("I9").Value = 1.044509891208315
$ws.Range("J9").Value = 1.06024815291053
$ws.Range("K9").Value = 1.05704209750602
$ws.Range("L9").Value = 1.063660503569154
$ws.Range("M9").Value = 1.07370823358823
$ws.Range("N9").Value = 1.061753826899826

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051858778027745
$ws.Range("D10").Value = 1.051839164384542
$ws.Range("E10").Value = 1.058124643669249
$ws.Range("F10").Value = 1.06804580209826
$ws.Range("I10").Value = 1.043701510364814
$ws.Range("J10").Value = 1.058083248497235
$ws.Range("K10").Value = 1.055229109904061
$ws.Range("L10").Value = 1.061493012283685
$ws.Range("M10").Value = 1.071380680801736
$ws.Range("N10").Value = 1.059585848073946

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050683641347929
$ws.Range("D11").Value = 1.050927011859769
$ws.Range("E11").Value = 1.057058935588193
$ws.Range("F11").Value = 1.066912251171927
$ws.Range("I11").Value = 1.043347751913728
$ws.Range("J11").Value = 1.057141254047852
$ws.Range("K11").Value = 1.054439646489491
$ws.Range("L11").Value = 1.060549686497537
$ws.Range("M11").Value = 1.070368417724803
$ws.Range("N11").Value = 1.058642515884397

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.05024647704506
$ws.Range("D12").Value = 1.050587613295075
$ws.Range("E12").Value = 1.056662432381907
$ws.Range("F12").Value = 1.066490600057788
$ws.Range("I12").Value = 1.043215784253414
$ws.Range("J12").Value = 1.056790652450909
$ws.Range("K12").Value = 1.054145726655572
$ws.Range("L12").Value = 1.060198558273894
$ws.Range("M12").Value = 1.069991738335236
$ws.Range("N12").Value = 1.058291416392972

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050340280725941
$ws.Range("D13").Value = 1.050660442108241
$ws.Range("E13").Value = 1.056747513432143
$ws.Range("F13").Value = 1.066581073077125
$ws.Range("I13").Value = 1.043244117483376
$ws.Range("J13").Value = 1.056865889740274
$ws.Range("K13").Value = 1.054208804343512
$ws.Range("L13").Value = 1.060273909985538
$ws.Range("M13").Value = 1.070072568367236
$ws.Range("N13").Value = 1.058366760527917

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050647518911182
$ws.Range("D14").Value = 1.050898969071066
$ws.Range("E14").Value = 1.057026173903558
$ws.Range("F14").Value = 1.066877409696644
$ws.Range("I14").Value = 1.043336855010304
$ws.Range("J14").Value = 1.057112287639157
$ws.Range("K14").Value = 1.054415364889
$ws.Range("L14").Value = 1.060520677214911
$ws.Range("M14").Value = 1.070337295251034
$ws.Range("N14").Value = 1.058613508340079

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050836729647046
$ws.Range("D15").Value = 1.051045855545397
$ws.Range("E15").Value = 1.057197778943542
$ws.Range("F15").Value = 1.067059912462995
$ws.Range("I15").Value = 1.043393918462288
$ws.Range("J15").Value = 1.057264007911524
$ws.Range("K15").Value = 1.054542543436856
$ws.Range("L15").Value = 1.060672620777936
$ws.Range("M15").Value = 1.07050031165308
$ws.Range("N15").Value = 1.058765444072636

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05193667559513
$ws.Range("D16").Value = 1.051899619721229
$ws.Range("E16").Value = 1.058195280730992
$ws.Range("F16").Value = 1.068120948978591
$ws.Range("I16").Value = 1.043724909130373
$ws.Range("J16").Value = 1.058145667838491
$ws.Range("K16").Value = 1.055281409617831
$ws.Range("L16").Value = 1.061555515499288
$ws.Range("M16").Value = 1.07144776692517
$ws.Range("N16").Value = 1.059648356057826

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052625475849819
$ws.Range("D17").Value = 1.052434136441291
$ws.Range("E17").Value = 1.058819843949809
$ws.Range("F17").Value = 1.068785458964999
$ws.Range("I17").Value = 1.043931529191852
$ws.Range("J17").Value = 1.05869747377145
$ws.Range("K17").Value = 1.055743686453734
$ws.Range("L17").Value = 1.062108039173946
$ws.Range("M17").Value = 1.072040886588444
$ws.Range("N17").Value = 1.060200945618503

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053026826973911
$ws.Range("D18").Value = 1.052745545232323
$ws.Range("E18").Value = 1.059183734488114
$ws.Range("F18").Value = 1.069172683494378
$ws.Range("I18").Value = 1.044051688289008
$ws.Range("J18").Value = 1.059018892310909
$ws.Range("K18").Value = 1.056012898359392
$ws.Range("L18").Value = 1.0624298560073
$ws.Range("M18").Value = 1.072386418142648
$ws.Range("N18").Value = 1.060522820609155

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053163607613387
$ws.Range("D19").Value = 1.052851665948476
$ws.Range("E19").Value = 1.059307743294424
$ws.Range("F19").Value = 1.069304654322788
$ws.Range("I19").Value = 1.044092598757492
$ws.Range("J19").Value = 1.059128413467279
$ws.Range("K19").Value = 1.056104620728596
$ws.Range("L19").Value = 1.06253950953089
$ws.Range("M19").Value = 1.072504163853021
$ws.Range("N19").Value = 1.06063249729813

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052551617102246
$ws.Range("D20").Value = 1.052376825786599
$ws.Range("E20").Value = 1.05875287640739
$ws.Range("F20").Value = 1.068714202028551
$ws.Range("I20").Value = 1.04390939798316
$ws.Range("J20").Value = 1.058638315893391
$ws.Range("K20").Value = 1.055694132715672
$ws.Range("L20").Value = 1.062048806396905
$ws.Range("M20").Value = 1.071977294549554
$ws.Range("N20").Value = 1.060141703729472

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050557063467279
$ws.Range("D21").Value = 1.050828745070242
$ws.Range("E21").Value = 1.056944133468488
$ws.Range("F21").Value = 1.066790162646643
$ws.Range("I21").Value = 1.043309561797592
$ws.Range("J21").Value = 1.057039749127276
$ws.Range("K21").Value = 1.054354556775155
$ws.Range("L21").Value = 1.060448030833806
$ws.Range("M21").Value = 1.070259358678699
$ws.Range("N21").Value = 1.058540866815191

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049299144371489
$ws.Range("D22").Value = 1.049852016703865
$ws.Range("E22").Value = 1.055803125547856
$ws.Range("F22").Value = 1.065576963563167
$ws.Range("I22").Value = 1.042929142880548
$ws.Range("J22").Value = 1.056030592166743
$ws.Range("K22").Value = 1.053508383258332
$ws.Range("L22").Value = 1.059437299155479
$ws.Range("M22").Value = 1.069175282891789
$ws.Range("N22").Value = 1.057530276736022

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049966363729034
$ws.Range("D23").Value = 1.050370124647591
$ws.Range("E23").Value = 1.056408359354015
$ws.Range("F23").Value = 1.066220438976569
$ws.Range("I23").Value = 1.043131123063618
$ws.Range("J23").Value = 1.056565956758227
$ws.Range("K23").Value = 1.053957332563516
$ws.Range("L23").Value = 1.059973516287582
$ws.Range("M23").Value = 1.069750351011335
$ws.Range("N23").Value = 1.058066401606633

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052584991971254
$ws.Range("D24").Value = 1.052402723134611
$ws.Range("E24").Value = 1.05878313741341
$ws.Range("F24").Value = 1.068746401122922
$ws.Range("I24").Value = 1.043919399232609
$ws.Range("J24").Value = 1.058665048153034
$ws.Range("K24").Value = 1.055716525234125
$ws.Range("L24").Value = 1.062075572563041
$ws.Range("M24").Value = 1.072006030367973
$ws.Range("N24").Value = 1.06016847395199

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055607573021488
$ws.Range("D25").Value = 1.054747138696951
$ws.Range("E25").Value = 1.06152304230054
$ws.Range("F25").Value = 1.071663071116791
$ws.Range("I25").Value = 1.044820093665531
$ws.Range("J25").Value = 1.06108369878534
$ws.Range("K25").Value = 1.057741289313943
$ws.Range("L25").Value = 1.064496866726026
$ws.Range("M25").Value = 1.074607000135871
$ws.Range("N25").Value = 1.062590559345617
